$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.327.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.811.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.599.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.327.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.301.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -12.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.768"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.723.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0983"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("E51").Value = "  -0.31%  "
